# Issue #12 - "add a reward" slide.
#
# A new content slide ("Add a reward" / "The reward will impact the next
# decision for SAI" + "Explain with a schema") is inserted right after the
# current slide 8 ("issues") and before the current slide 9 ("Solutions"),
# i.e. it becomes the new slide 9 - "Solutions" and every slide after it
# simply shift one position later.

$p = $ppt.ActivePresentation

# ppLayoutText (2) -> "Titre et contenu" (title + content placeholder),
# the same layout used by every other content slide in this deck.
$s = $p.Slides.Add(9, 2)

$title = $s.Shapes.Placeholders.Item(1)
$title.Name = "Titre 1"
$title.TextFrame.TextRange.Text = "Add a reward"
$title.TextFrame.TextRange.LanguageID = "fr-FR"

$body = $s.Shapes.Placeholders.Item(2)
$body.Name = "Espace réservé du contenu 2"
$body.TextFrame.TextRange.Text = "The reward will impact the next decision for SAI`rExplain with a schema"
$body.TextFrame.TextRange.LanguageID = "fr-FR"
